# Angptl1-Tek LR-pair sheet: rebuilt against the new TPM run.
# - "Inflammatory-Mac" sending/target cluster is renamed "Resolving-Mac".
# - Every numeric column (E:T) is recomputed from the new TPM values.
# - The old 16-row table (Sending clusters ECs/FAPs/Inflammatory-Mac/MuSCs) becomes a
#   15-row table (Sending clusters FAPs/MuSCs/Resolving-Mac), so the trailing row is removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Delete()

$rows = @(
    ("FAPs", "Angptl1", "Tek", "ECs", 3, 1, 31.36337633333333, 94.090129, 0.9925580820795565, 0.9925580820795566, 3, 1, 58.51417433333334, 175.542523, 0.9208013000516164, 0.9208013000516164, 1835.202070450608, 16516.81863405547, 0.9139487723555947, 0.9139487723555948),
    ("FAPs", "Angptl1", "Tek", "FAPs", 3, 1, 31.36337633333333, 94.090129, 0.9925580820795565, 0.9925580820795566, 3, 1, 4.208408333333334, 12.625225, 0.06622511397676659, 0.06622511397676657, 131.9898943226695, 1187.909048904025, 0.06573227211427948, 0.06573227211427947),
    ("FAPs", "Angptl1", "Tek", "MuSCs", 3, 1, 31.36337633333333, 94.090129, 0.9925580820795565, 0.9925580820795566, 2, 0.6666666666666666, 0.1957283333333334, 0.5871850000000001, 0.003080055488155473, 0.003080055488155473, 6.138701377429446, 55.24831239686501, 0.003057133968022209, 0.003057133968022209),
    ("FAPs", "Angptl1", "Tek", "Neutrophils", 3, 1, 31.36337633333333, 94.090129, 0.9925580820795565, 0.9925580820795566, 3, 1, 0.6056613333333333, 1.816984, 0.00953091707228673, 0.00953091707228673, 18.99558432788178, 170.960258950936, 0.00945998876972822, 0.00945998876972822),
    ("FAPs", "Angptl1", "Tek", "Resolving-Mac", 3, 1, 31.36337633333333, 94.090129, 0.9925580820795565, 0.9925580820795566, 1, 0.3333333333333333, 0.023043, 0.069129, 0.0003626134111748421, 0.0003626134111748421, 0.7227062808489999, 6.504356527641, 0.0003599148719320269, 0.000359914871932027),
    ("MuSCs", "Angptl1", "Tek", "ECs", 2, 0.6666666666666666, 0.1453933333333333, 0.43618, 0.004601268898690329, 0.004601268898690329, 3, 1, 58.51417433333334, 175.542523, 0.9208013000516164, 0.9208013000516164, 8.507570853571112, 76.56813768214, 0.004236854383801125, 0.004236854383801125),
    ("MuSCs", "Angptl1", "Tek", "FAPs", 2, 0.6666666666666666, 0.1453933333333333, 0.43618, 0.004601268898690329, 0.004601268898690329, 3, 1, 4.208408333333334, 12.625225, 0.06622511397676659, 0.06622511397676657, 0.6118745156111113, 5.506870640500001, 0.0003047195572535183, 0.0003047195572535182),
    ("MuSCs", "Angptl1", "Tek", "MuSCs", 2, 0.6666666666666666, 0.1453933333333333, 0.43618, 0.004601268898690329, 0.004601268898690329, 2, 0.6666666666666666, 0.1957283333333334, 0.5871850000000001, 0.003080055488155473, 0.003080055488155473, 0.02845759481111112, 0.2561183533, 0.00001417216352389024, 0.00001417216352389024),
    ("MuSCs", "Angptl1", "Tek", "Neutrophils", 2, 0.6666666666666666, 0.1453933333333333, 0.43618, 0.004601268898690329, 0.004601268898690329, 3, 1, 0.6056613333333333, 1.816984, 0.00953091707228673, 0.00953091707228673, 0.08805912012444445, 0.79253208112, 0.00004385431230070962, 0.00004385431230070962),
    ("MuSCs", "Angptl1", "Tek", "Resolving-Mac", 2, 0.6666666666666666, 0.1453933333333333, 0.43618, 0.004601268898690329, 0.004601268898690329, 1, 0.3333333333333333, 0.023043, 0.069129, 0.0003626134111748421, 0.0003626134111748421, 0.00335029858, 0.03015268722, 0.000001668481811086809, 0.000001668481811086809),
    ("Resolving-Mac", "Angptl1", "Tek", "ECs", 2, 0.6666666666666666, 0.08976033333333333, 0.269281, 0.002840649021753016, 0.002840649021753016, 3, 1, 58.51417433333334, 175.542523, 0.9208013000516164, 0.9208013000516164, 5.252251792884778, 47.27026613596301, 0.00261567331222053, 0.00261567331222053),
    ("Resolving-Mac", "Angptl1", "Tek", "FAPs", 2, 0.6666666666666666, 0.08976033333333333, 0.269281, 0.002840649021753016, 0.002840649021753016, 3, 1, 4.208408333333334, 12.625225, 0.06622511397676659, 0.06622511397676657, 0.3777481348027778, 3.399733213225, 0.000188122305233584, 0.0001881223052335839),
    ("Resolving-Mac", "Angptl1", "Tek", "MuSCs", 2, 0.6666666666666666, 0.08976033333333333, 0.269281, 0.002840649021753016, 0.002840649021753016, 2, 0.6666666666666666, 0.1957283333333334, 0.5871850000000001, 0.003080055488155473, 0.003080055488155473, 0.01756864044277778, 0.158117763985, 0.000008749356609373853, 0.000008749356609373852),
    ("Resolving-Mac", "Angptl1", "Tek", "Neutrophils", 2, 0.6666666666666666, 0.08976033333333333, 0.269281, 0.002840649021753016, 0.002840649021753016, 3, 1, 0.6056613333333333, 1.816984, 0.00953091707228673, 0.00953091707228673, 0.0543643631671111, 0.4892792685039999, 0.00002707399025780042, 0.00002707399025780042),
    ("Resolving-Mac", "Angptl1", "Tek", "Resolving-Mac", 2, 0.6666666666666666, 0.08976033333333333, 0.269281, 0.002840649021753016, 0.002840649021753016, 1, 0.3333333333333333, 0.023043, 0.069129, 0.0003626134111748421, 0.0003626134111748421, 0.002068347361, 0.018615126249, 0.000001030057431728339, 0.000001030057431728339)
)

$data = New-Object 'object[,]' $rows.Count,20
for ($r = 0; $r -lt $rows.Count; $r++) {
    for ($c = 0; $c -lt 20; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

$lastRow = 1 + $rows.Count
$ws.Range("A2:T$lastRow").Value = $data
